# "End of the day"
#  - the last existing task now finishes a bit later (11:20 instead of 11:15)
#  - a new "Video" task is appended, running 11:20 -> 11:40
#  - the whole table is left aligned (time/description columns vertically
#    centered too), and the header row keeps its bold/bordered look but is
#    also left aligned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft         = -4131   # xlHAlignLeft
$xlCenterV      = -4108   # xlVAlignCenter
$xlPasteFormats = -4122   # xlPasteFormats

# ---------------------------------------------------------------------
# Build each distinct target format once on an out-of-the-way scratch
# cell (touched exactly once each, so no stray in-between styles get
# created), then stamp the finished format onto the real cells with a
# single "paste formats" operation per destination range.
# ---------------------------------------------------------------------

# Scratch cell Z1 -> time format, left+vcenter aligned (rows 3-8, cols B/C)
$ws.Range("Z1").NumberFormat = "h:mm"
$ws.Range("Z1").HorizontalAlignment = $xlLeft
$ws.Range("Z1").VerticalAlignment = $xlCenterV
$ws.Range("Z1").Copy()
$ws.Range("B3:C8").PasteSpecial($xlPasteFormats)

# Scratch cell Z2 -> general format, left+vcenter aligned (description column)
$ws.Range("Z2").HorizontalAlignment = $xlLeft
$ws.Range("Z2").VerticalAlignment = $xlCenterV
$ws.Range("Z2").Copy()
$ws.Range("D3:D8").PasteSpecial($xlPasteFormats)

# Header cells only need left alignment added on top of their existing
# bold font + border, so a single property write is enough for each.
$ws.Range("B2").HorizontalAlignment = $xlLeft
$ws.Range("C2").HorizontalAlignment = $xlLeft
$ws.Range("D2").HorizontalAlignment = $xlLeft

# Scratch cell Z3 -> time format, left aligned only (new row 9, cols B/C)
$ws.Range("Z3").NumberFormat = "h:mm"
$ws.Range("Z3").HorizontalAlignment = $xlLeft

# ---------------------------------------------------------------------
# Data edits: push row 8's end time out, and add the new "Video" row.
# ---------------------------------------------------------------------
$ws.Range("C8").Value = 0.47222222222222227

$ws.Range("B9").Value = 0.47222222222222227
$ws.Range("C9").Value = 0.4861111111111111
$ws.Range("D9").Value = "Video"

$ws.Range("Z3").Copy()
$ws.Range("B9:C9").PasteSpecial($xlPasteFormats)

$ws.Range("Z2").Copy()
$ws.Range("D9").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Clean up the scratch cells so they don't show up in the saved sheet.
# ---------------------------------------------------------------------
$ws.Range("Z1").Clear()
$ws.Range("Z2").Clear()
$ws.Range("Z3").Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# View state: selection moved to H4, sheet zoomed to 160%.
# ---------------------------------------------------------------------
$ws.Range("H4").Select()
$excel.ActiveWindow.Zoom = 160
